$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Пример")

# Update the organization website from "www.stat.kg" to "www.stat.gov.kg"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Reflect the active cell selection change observed in the saved file
$ws.Range("B7").Select()
